$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.561.71"
$ws.Range("E2").Value = "  -2.05%  "
$ws.Range("D3").Value = "2.971.16"
$ws.Range("E3").Value = "  -2.24%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "587.59"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.03%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "141.31"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -6.02%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("E8").Value = "  -2.07%  "
$ws.Range("D9").Value = "2.966.76"
$ws.Range("E9").Value = "  -2.45%  "
$ws.Range("E10").Value = "  -5.76%  "
$ws.Range("E11").Value = "  -0.35%  "
$ws.Range("E12").Value = "  +3.00%  "
$ws.Range("E13").Value = "  -2.65%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "33.96"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -5.08%  "
$ws.Range("E15").Value = "  +1.42%  "
$ws.Range("D16").Value = "3.467.54"
$ws.Range("E16").Value = "  -2.15%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "7.00"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.21%  "
$ws.Range("D18").Value = "61.567.78"
$ws.Range("E18").Value = "  -2.09%  "
$ws.Range("D19").Value = "2.976.18"
$ws.Range("E19").Value = "  -2.43%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "449.17"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -5.83%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.87"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.32%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.682"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.71%  "
$ws.Range("E23").Value = "  -2.02%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "81.18"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.45%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.11"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.49%  "
$ws.Range("E26").Value = "  -9.51%  "
$ws.Range("E27").Value = "  +0.20%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.80"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -6.69%  "
$ws.Range("E29").Value = "  -0.18%  "
$ws.Range("E30").Value = "  -0.31%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.83"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -6.64%  "
$ws.Range("E32").Value = "  -5.62%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "27.10"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.66%  "
$ws.Range("E34").Value = "  -2.64%  "
$ws.Range("E35").Value = "  -4.21%  "
$ws.Range("D36").Value = "0.0₃0773"
$ws.Range("E36").Value = "  -3.67%  "
$ws.Range("E37").Value = "  -2.78%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.07"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.68%  "
$ws.Range("E39").Value = "  +0.34%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "50.08"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.36%  "
$ws.Range("B41").Value = "Kaspa"
$ws.Range("C41").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.119"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.57%  "
$ws.Range("B42").Value = "dogwifhat"
$ws.Range("C42").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.79"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -10.37%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "387.56"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -8.60%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0353"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.31%  "
$ws.Range("D45").Value = "2.725.37"
$ws.Range("E45").Value = "  -3.88%  "
$ws.Range("E46").Value = "  -7.56%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "36.91"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.68%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "129.85"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.65%  "
$ws.Range("E49").Value = "  +0.06%  "
$ws.Range("E50").Value = "  -1.18%  "
$ws.Range("E51").Value = "  -0.43%  "
